# Update the "想去人数" (interested-count) figures in column F for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets to the freshly
# scraped values, as published by the gh-pages data-refresh commit.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F, shared by both worksheets
# (the underlying events are identical; only row offsets differ because
# the "全部类型" sheet contains two additional rows).
$updates1 = @{
    2  = 240
    3  = 539
    4  = 13996
    5  = 236
    6  = 1806
    11 = 556
    12 = 37
    15 = 14136
    16 = 377
    18 = 15023
    19 = 18
    20 = 8370
    21 = 285
    24 = 159
    26 = 169
    28 = 14
    30 = 36
    31 = 1046
    32 = 27
    35 = 408
    37 = 13
    39 = 232
    40 = 398
    42 = 5152
}

$updates4 = @{
    2  = 240
    3  = 539
    4  = 13996
    5  = 236
    6  = 1806
    11 = 556
    12 = 37
    15 = 14136
    16 = 377
    18 = 15023
    19 = 18
    20 = 8370
    21 = 285
    24 = 159
    26 = 169
    28 = 14
    30 = 36
    31 = 1046
    32 = 27
    37 = 408
    39 = 13
    41 = 232
    42 = 398
    44 = 5152
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
